# Updates cryptos list prices/volumes (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells below contain plain numeric-looking text (e.g. "1.00", "0.0500")
# that Excel would otherwise coerce to a Number and mangle (dropping trailing
# zeros / introducing float noise). Mark them Text first so the literal string
# is preserved exactly, matching the source data (which is always text).
$textCells = @("D4", "D5", "D6", "D10", "D13", "D15", "D17", "D20", "D22", "D26", "D27", "D29", "D31", "D34", "D38", "D42", "D46", "D49")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.006.92'
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").Value = '1.642.50'
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.70%  '
$ws.Range("D5").Value = '216.26'
$ws.Range("D6").Value = '0.507'
$ws.Range("E6").Value = '  +1.03%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("E9").Value = '  +0.84%  '
$ws.Range("D10").Value = '19.64'
$ws.Range("E10").Value = '  -0.41%  '
$ws.Range("E11").Value = '  +0.78%  '
$ws.Range("B12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D12").Value = '1.871.36'
$ws.Range("E12").Value = '  +0.24%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '4.29'
$ws.Range("E13").Value = '  +0.44%  '
$ws.Range("D14").Value = '1.676.39'
$ws.Range("E14").Value = '  +4.48%  '
$ws.Range("D15").Value = '0.545'
$ws.Range("E15").Value = '  +0.38%  '
$ws.Range("D16").Value = '0.0₃0766'
$ws.Range("E16").Value = '  +0.75%  '
$ws.Range("D17").Value = '63.02'
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("D18").Value = '25.949.46'
$ws.Range("E18").Value = '  -0.50%  '
$ws.Range("E19").Value = '  +0.66%  '
$ws.Range("D20").Value = '193.11'
$ws.Range("E20").Value = '  -0.78%  '
$ws.Range("E21").Value = '  -1.26%  '
$ws.Range("D22").Value = '9.93'
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("E24").Value = '  +6.77%  '
$ws.Range("E25").Value = '  +1.21%  '
$ws.Range("D26").Value = '144.60'
$ws.Range("E26").Value = '  +1.44%  '
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.59%  '
$ws.Range("E28").Value = '  +0.65%  '
$ws.Range("D29").Value = '15.56'
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  +0.45%  '
$ws.Range("D31").Value = '0.0500'
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("E32").Value = '  -1.55%  '
$ws.Range("E33").Value = '  +0.65%  '
$ws.Range("D34").Value = '1.53'
$ws.Range("E34").Value = '  -3.04%  '
$ws.Range("E35").Value = '  +2.33%  '
$ws.Range("E36").Value = '  -0.63%  '
$ws.Range("D37").Value = '1.134.81'
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("D38").Value = '0.543'
$ws.Range("E38").Value = '  -1.37%  '
$ws.Range("E39").Value = '  -0.73%  '
$ws.Range("E40").Value = '  +0.25%  '
$ws.Range("E41").Value = '  +0.68%  '
$ws.Range("D42").Value = '99.69'
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("D44").Value = '1.780.29'
$ws.Range("E44").Value = '  +0.20%  '
$ws.Range("E45").Value = '  +2.72%  '
$ws.Range("D46").Value = '56.75'
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("E47").Value = '  +2.60%  '
$ws.Range("E48").Value = '  -0.27%  '
$ws.Range("D49").Value = '7.75'
$ws.Range("E49").Value = '  +1.13%  '
$ws.Range("E50").Value = '  -0.10%  '
$ws.Range("E51").Value = '  -0.08%  '

# Restore default (unstyled) cell style now that the text values are locked in,
# so formatting matches the rest of the untouched sheet.
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
